$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data rows down
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the row below; strip it so the
# new row matches the plain (unstyled) look of the other data rows.
$ws.Rows("2:2").ClearFormats()

# Keep the date as literal text (matches existing rows) instead of letting
# Excel auto-convert the "yyyy-mm-dd" string into a date serial number.
$ws.Range("A2").NumberFormat = "@"

# Fill the new row with the latest price entry (one day after the prior top row)
$ws.Range("A2").Value = "2025-12-28"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610

# Drop back to the default/general format now that the text value is locked in,
# so the new row has no leftover style (matching the other plain data rows).
$ws.Range("A2").ClearFormats()
